# "editing geneSearch for both common and rare inversion breakpoints"
#
# The "Unique set" block at the bottom of the sheet (rows 53-58) used to be
# separated from the table above by a completely blank row 53, and its own
# header row (54) carried formatted-but-empty cells in B:I. This edit drops
# that blank separator row - shifting the whole "Unique set" block up by one
# row (54->53 ... 58->57) - and strips the now-pointless empty formatted
# cells from the header row, leaving it with only its label in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 is entirely empty in the source sheet - it only exists as a visual
# gap above the "Unique set" header. Deleting it shifts rows 54:58 up to
# 53:57 (values, number/text types, row heights and formatting all move
# with it), which reproduces the whole block exactly as the diff wants it.
$ws.Rows.Item(53).Delete()

# The header row ("Unique set", now row 53) no longer needs its empty,
# formatted B:I filler cells - clear them (content + formatting) so the row
# only carries the column-A label, matching the new layout.
$ws.Range("B53:I53").Clear()

# Reflect the saved selection/active cell from the edit.
$ws.Range("I19").Select() | Out-Null
